$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 429, shifting existing rows 429:490 down to 430:491.
$ws.Rows("429:429").Insert()

# Populate the newly inserted row 429 with the new record's data.
$ws.Cells.Item(429, 1).Value = 3
$ws.Cells.Item(429, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(429, 3).Value = "Coquimbo"
$ws.Cells.Item(429, 4).Value = 45127
$ws.Cells.Item(429, 5).Value = 5
$ws.Cells.Item(429, 6).Value = 100112001
$ws.Cells.Item(429, 7).Value = "Berenjena"
$ws.Cells.Item(429, 8).Value = "Sin especificar"
$ws.Cells.Item(429, 9).Value = "Primera"
$ws.Cells.Item(429, 10).Value = 50
$ws.Cells.Item(429, 11).Value = 7000
$ws.Cells.Item(429, 12).Value = 7000
$ws.Cells.Item(429, 13).Value = 7000
$ws.Cells.Item(429, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(429, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(429, 16).Value = 117
$ws.Cells.Item(429, 17).Value = 60
$ws.Cells.Item(429, 18).Value = "Hortaliza"
